# "Update 13 Jul 2021, end of day update."
# Petty cash book - Sheet1 ("Buku Kas Harian" main sheet) gets a day's
# worth of new transactions plus corrections to two already-entered rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: the 60,000 cash entry actually also included a 240,000 item ---
$ws.Range("D3").Formula = "=60000+240000"

# --- Row 6: more debit items recorded for that running total ---
$ws.Range("D6").Formula = "=406000+280000+4600000+5000000+3000000+41430000+3490000+580000+6560000+6031000+6027000+1040000"

# --- Row 7: more credit items recorded for that running total ---
$ws.Range("C7").Formula = "=5000000+875000+41430000+6560000+6031000+6027000+37540000"

# --- Row 9: new transaction - PRIVE (owner's draw) by andreas ---
$ws.Range("B9").Value = "PRIVE - andreas"
$ws.Range("D9").Value = 5000000

# --- Row 10: new transaction - cash/retail sales ---
$ws.Range("B10").Value = "SALES - cash/retail"
$ws.Range("C10").Formula = "=22097025+18794975-37540000"

# --- Row 11: new transaction - cash overage (selisih lebih) ---
$ws.Range("B11").Value = "SELISIH - lebih"
$ws.Range("C11").Value = 30000

# --- Row 12: new transaction - deposit to bank ---
$ws.Range("B12").Value = "SETOR KE BANK"
$ws.Range("D12").Value = 21000000

# --- Row 13: new day starts - 13 Jul 2021 ---
$ws.Range("A13").Value = 44390

# Leave the cursor where the day's last manual entry was made.
$ws.Range("C11").Select()
